$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'44.624.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.85%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.244.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.04%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.33%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'306.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.34%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'94.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.45%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.49%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.10%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.515"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.28%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'34.89"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.13%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0800"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.47%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'7.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.34%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -0.02%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.587.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.02%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.245.48"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.85%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.831"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.19%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'13.54"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.18%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'44.409.89"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.96%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -3.02%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -3.46%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'11.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -4.02%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'65.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.52%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'237.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.67%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.31%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'1.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.91%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.19%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +4.59%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'9.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.15%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'36.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -3.79%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'19.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.75%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'5.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.36%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'147.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -4.17%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.0782"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -2.00%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -0.04%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'3.19"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.29%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.108"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.06%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -1.84%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +4.70%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'15.18"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +4.79%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -5.15%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'3.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.96%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.0299"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.37%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +0.11%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.812.39"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +3.55%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +11.97%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'81.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.39%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.187"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.33%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'98.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.05%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'68.88"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +2.25%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'4.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.11%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'54.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.23%  "
$ws.Range("E51").Style = "Normal"
